# Update "MCF" maximum capacity factor values to 1 and bump the "About"
# sheet's last-updated date. Mirrors the workbook author's 4.0 refresh.

$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date in C1 ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45392

# --- "MCF" sheet: set capacity-factor inputs to 1 ---
$mcf = $wb.Worksheets.Item("MCF")

$cells = @("B2", "B3", "B4", "B6", "B10", "B11", "B12", "B13", "B14", "B16", "B17", "B18")
foreach ($cellRef in $cells) {
    $mcf.Range($cellRef).Value = 1
}

# Formula-driven cells (B19, B20, B21, B22, B24, B25) recalc automatically
# from the edits above since they reference B2/B4/B10/B14.

# Restore the active sheet / selection to match the saved workbook state.
$mcf.Activate()
$mcf.Range("B17").Select()
